$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3: B3 = SUM(B1:B2)
$ws.Range("B3").Formula = "=SUM(B1:B2)"

# Update the selection shown in the sheet view
$ws.Range("G5").Select()
